$wb = $excel.ActiveWorkbook

# --- Weekly Quantity: append 4 new weekly rows (41-44) ---
$ws = $wb.Worksheets.Item("Weekly Quantity")
$weeklyNewRows = @(
    @(45662.99999999999, 12),
    @(45669.99999999999, 7),
    @(45676.99999999999, 1),
    @(45683.99999999999, 2)
)
$startRow = 41
for ($i = 0; $i -lt $weeklyNewRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $weeklyNewRows[$i][0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $weeklyNewRows[$i][1]
}

# --- Monthly Trend: append 1 new monthly row (18) ---
$ws = $wb.Worksheets.Item("Monthly Trend")
$ws.Cells.Item(18, 1).Value = 45688.99999999999
$ws.Cells.Item(18, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 2).Value = 22

# --- PO Forecast: new forecast model -- full column rewrite + 4 new rows (49-52) ---
$ws = $wb.Worksheets.Item("PO Forecast")
$forecastRows = @(
    @(45151.99999999999, 18),
    @(45186.99999999999, 18),
    @(45193.99999999999, 18),
    @(45207.99999999999, 18),
    @(45214.99999999999, 18),
    @(45221.99999999999, 18),
    @(45228.99999999999, 18),
    @(45235.99999999999, 18),
    @(45242.99999999999, 18),
    @(45270.99999999999, 18),
    @(45277.99999999999, 18),
    @(45298.99999999999, 18),
    @(45305.99999999999, 18),
    @(45312.99999999999, 18),
    @(45319.99999999999, 18),
    @(45326.99999999999, 18),
    @(45333.99999999999, 18),
    @(45340.99999999999, 18),
    @(45347.99999999999, 18),
    @(45368.99999999999, 18),
    @(45375.99999999999, 18),
    @(45382.99999999999, 18),
    @(45389.99999999999, 18),
    @(45459.99999999999, 18),
    @(45487.99999999999, 18),
    @(45501.99999999999, 18),
    @(45515.99999999999, 18),
    @(45529.99999999999, 18),
    @(45536.99999999999, 18),
    @(45543.99999999999, 18),
    @(45564.99999999999, 17),
    @(45571.99999999999, 17),
    @(45578.99999999999, 17),
    @(45585.99999999999, 17),
    @(45599.99999999999, 17),
    @(45606.99999999999, 17),
    @(45613.99999999999, 17),
    @(45634.99999999999, 17),
    @(45641.99999999999, 17),
    @(45662.99999999999, 17),
    @(45669.99999999999, 17),
    @(45676.99999999999, 17),
    @(45683.99999999999, 17),
    @(45690.99999999999, 17),
    @(45697.99999999999, 17),
    @(45704.99999999999, 17),
    @(45711.99999999999, 17),
    @(45718.99999999999, 17),
    @(45725.99999999999, 17),
    @(45732.99999999999, 17),
    @(45739.99999999999, 17)
)
for ($i = 0; $i -lt $forecastRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $forecastRows[$i][0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $forecastRows[$i][1]
}

